# Appends the 21 new purchase-data rows (925-945) that were added to
# Sheet1 in the "Add files via upload" commit: four new Ref/part
# numbers (GNA80496, GNA80489, GNA80487, GNA80488), each with their own
# Date / Qty / Unit Price / Total line items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    # Ref,        Date,    Qty, UnitPrice,            Total
    @("GNA80496", 45925,  15,  92,                    1380),
    @("GNA80496", 45925,  15,  215,                   3225),
    @("GNA80496", 45925,  50,  114,                   5700),
    @("GNA80496", 45925,  45,  132,                   5940),
    @("GNA80496", 45925,  60,  161,                   9660),
    @("GNA80496", 45925,  45,  198,                   8910),
    @("GNA80496", 45925,  10,  241,                   2410),
    @("GNA80496", 45925,  30,  93,                    2790),
    @("GNA80489", 45923,  120, 11.5,                  1380),
    @("GNA80489", 45923,  300, 15.24,                 4572),
    @("GNA80489", 45923,  150, 19,                    2850),
    @("GNA80489", 45923,  220, 21.03,                 4626.6000000000004),
    @("GNA80489", 45923,  25,  21.71,                 542.75),
    @("GNA80487", 45923,  50,  21.55,                 1077.5),
    @("GNA80487", 45923,  50,  23,                    1150),
    @("GNA80487", 45923,  50,  25.42,                 1271),
    @("GNA80487", 45923,  50,  19.829999999999998,    991.5),
    @("GNA80488", 45923,  130, 21.55,                 2801.5),
    @("GNA80488", 45923,  50,  23,                    1150),
    @("GNA80488", 45923,  98,  25.42,                 2491.16),
    @("GNA80488", 45923,  50,  19.829999999999998,    1983)
)

$startRow = 925
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Match the author's final selection/scroll position from the saved file.
$lastRow = $startRow + $newRows.Count - 1
$ws.Range("E$lastRow").Select() | Out-Null
